$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. "Színek" header (row 6) gets the same header style as "Főoldal
#    dizájnja" (row 2) / the new "Betűtípusok" header (row 11):
#    bigger font + bottom border (style index 1 in styles.xml).
#    Copying the format from B2 reuses the existing style instead of
#    creating a duplicate one.
# ------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Rows(6).RowHeight = 18

# ------------------------------------------------------------------
# 2. New "Betűtípusok" (fonts) section, starting at row 11.
# ------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Rows(11).RowHeight = 18
$ws.Range("B11").Value = "Betűtípusok"

# Header row (row 12) - reuses B2's bottom border but keeps the
# normal (non-enlarged) font, producing the new style index 5.
$ws.Range("B2").Copy()
$ws.Range("B12:E12").PasteSpecial(-4122)
$ws.Range("B12:E12").Font.Size = 11

# Fill in the table values in the same order the workbook's author
# originally typed them (this controls shared-string ordering).
$ws.Range("B12").Value = "Elem"
$ws.Range("C12").Value = "Betűtípus"
$ws.Range("D12").Value = "Betűstílus"

$ws.Range("C13").Value = "Poppins"
$ws.Range("D14").Value = "Regular 400"
$ws.Range("B13").Value = "h1, h2, h3 címek"
$ws.Range("D13").Value = "Regular 500"
$ws.Range("B14").Value = "p (és minden egyéb szöveg)"
$ws.Range("B15").Value = "navbar szövege"
$ws.Range("C15").Value = "Open Sans"
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("C14").Value = "Raleway"
$ws.Range("E12").Value = "Fonts link"
$ws.Range("E14").Value = "https://fonts.google.com/specimen/Raleway"
$ws.Range("E13").Value = "https://fonts.google.com/specimen/Poppins"

# ------------------------------------------------------------------
# 3. Column widths (best-fit) for the new content.
# ------------------------------------------------------------------
$ws.Columns(2).AutoFit()
$ws.Columns(4).AutoFit()

# ------------------------------------------------------------------
# 4. Leave the selection where the author last left it.
# ------------------------------------------------------------------
[void]$ws.Range("F10").Select()
